# Build site at 2022-09-26 16:07:08 UTC
# Restructure the LOM3085 "Objetivos / Docentes responsaveis / Programa / Avaliacao / Bibliografia"
# block: drop the long narrative paragraphs, and shift the four professor-name
# rows up so each one sits alongside its section label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Drop the five rows that disappear entirely (old rows 23-27), shrinking
#    the sheet from A1:C27 down to A1:C22.
# ---------------------------------------------------------------------------
$ws.Range("A23:C27").EntireRow.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2) Cells that exist now but must not exist in the final layout: remove them
#    completely (not just their contents) so no stray empty <c> survives.
# ---------------------------------------------------------------------------
$ws.Range("B13").Clear() | Out-Null
$ws.Range("C13").Clear() | Out-Null
$ws.Range("B15").Clear() | Out-Null
$ws.Range("C15").Clear() | Out-Null
$ws.Range("B16").Clear() | Out-Null
$ws.Range("C16").Clear() | Out-Null
$ws.Range("A22").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 3) Cells that need to be created fresh: copy the correct column's
#    formatting in first (one cell at a time - multi-area PasteSpecial
#    targets are unreliable) so they don't inherit the wrong default style,
#    then they'll be given their text in step 4.
# ---------------------------------------------------------------------------
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B12").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C12").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("A17").Copy() | Out-Null
$ws.Range("A13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A16").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B10").Copy() | Out-Null
$ws.Range("B18").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C18").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B10").Copy() | Out-Null
$ws.Range("B20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C20").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Write the final text for rows 10-22.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("A11").Value = "Objectives:"

$ws.Range("A12").Value = "Programa resumido:"
$ws.Range("B12").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("C12").Value = "5840963 - Daniela Camargo Vernilli"

$ws.Range("A13").Value = "Short syllabus:"

$ws.Range("A14").Value = "Programa:"
$ws.Range("B14").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C14").Value = "6495737 - Durval Rodrigues Junior"

$ws.Range("A15").Value = "Syllabus:"

$ws.Range("A16").Value = "Avaliação:"

$ws.Range("A17").Value = "Método:"
$ws.Range("B17").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C17").Value = "984972 - Hugo Ricardo Zschommler Sandim"

$ws.Range("A18").Value = "Critério:"
$ws.Range("B18").Value = "Aulas expositivas complementadas com experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento e de estudo de casos."
$ws.Range("C18").Value = "Aulas expositivas complementadas com experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento e de estudo de casos."

$ws.Range("A19").Value = "Norma de recuperação:"
$ws.Range("B19").Value = "Média aritmética das notas obtidas nos relatórios e trabalhos. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."
$ws.Range("C19").Value = "Média aritmética das notas obtidas nos relatórios e trabalhos. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."

$ws.Range("A20").Value = "Bibliografia:"
$ws.Range("B20").Value = "Devido às características práticas da disciplina, não será oferecida recuperação."
$ws.Range("C20").Value = "Devido às características práticas da disciplina, não será oferecida recuperação."

$ws.Range("A21").Value = "Requisitos:"

$ws.Range("B22").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Range("C22").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"

# ---------------------------------------------------------------------------
# 5) Row heights: only a handful actually change value versus the original.
# ---------------------------------------------------------------------------
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 120
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 30

Write-Output "edit applied"
